# Edit: split the "Umrechnung des Abstands von „Pixel“ auf tatsächlicher Meter"
# paragraph into three runs, correcting "tatsächlicher" -> "tatsächliche":
#   Run 1: "Umrechnung des Abstands von „Pixel“ "
#   Run 2: "auf tatsächliche "
#   Run 3: "Meter"

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -like "*tatsächlicher Meter*") {
                $targetShape = $shape
                $targetSlide = $slide
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the phrase that needs to be rewritten/split.
$oldPhrase = "auf tatsächlicher Meter"
$newPhrase = "auf tatsächliche Meter"

$found = $tr.Find($oldPhrase, 0)
$start = $found.Start
$len = $found.Length

# Rewrite "auf tatsächlicher Meter" -> "auf tatsächliche Meter".
# This splits the original single run into two runs: the unchanged prefix
# ("Umrechnung des Abstands von „Pixel“ ") and the rewritten remainder
# ("auf tatsächliche Meter"), both keeping the original run formatting.
$sub = $tr.Characters($start, $len)
$sub.Text = $newPhrase

# Now split "Meter" out of the rewritten remainder into its own run, so the
# final structure is three runs:
#   "Umrechnung des Abstands von „Pixel“ " | "auf tatsächliche " | "Meter"
$meterStart = $start + $newPhrase.Length - 5
$meterPart = $tr.Characters($meterStart, 5)
$meterPart.Text = "Meter"
